# Quarterly income-statement "rolling window" update:
# drop the oldest quarter column, shift the remaining 9 quarters one
# column to the left, and append the newest quarter (1402-02-30) in
# column M -- for both the header rows (period name / publish date)
# and every numeric data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("D", "E", "F", "G", "H", "I", "J", "K", "L", "M")

# Row 8: "دوره مالی" (financial period) labels
$periodLabels = @(
    "فصل سوم منتهی به 1399/09",
    "فصل چهارم منتهی به 1399/12",
    "فصل اول منتهی به 1400/03",
    "فصل دوم منتهی به 1400/06",
    "فصل سوم منتهی به 1400/09",
    "فصل چهارم منتهی به 1400/12",
    "فصل اول منتهی به 1401/03",
    "فصل دوم منتهی به 1401/06",
    "فصل سوم منتهی به 1401/09",
    "فصل چهارم منتهی به 1401/12"
)

# Row 9: "تاریخ انتشار" (publish date) labels
$publishDates = @(
    "1400-10-29 (3)",
    "1401-03-11 (8)",
    "1401-04-29 (2)",
    "1401-08-29 (4)",
    "1401-10-28 (2)",
    "1402-02-30 (7)",
    "1401-04-29",
    "1401-08-29 (2)",
    "1401-10-28",
    "1402-02-30"
)

for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "8").Value = $periodLabels[$i]
    $ws.Range($cols[$i] + "9").Value = $publishDates[$i]
}

# Numeric / dash data rows 11-27, each already shifted one column left
# with the new (tenth) quarter value appended at the end (column M).
$rowData = @{
    11 = @(9363, 13836, 10811, 12840, 13057, 14887, 14707, 20016, 19121, 16441)
    12 = @(-8277, -10470, -9171, -10582, -11332, -11587, -13082, -13748, -16376, -13173)
    13 = @(1086, 3366, 1640, 2258, 1725, 3300, 1625, 6268, 2745, 3268)
    14 = @(-886, -1525, -841, -1054, -915, -1312, -875, -1350, -884, -829)
    15 = @("-", "-", "-", "-", "-", "-", "-", "-", "-", "-")
    16 = @(339, -191, -584, 533, 113, -110, 77, 43, 772, 1274)
    17 = @(539, 1650, 215, 1737, 922, 1878, 827, 4961, 2632, 3713)
    18 = @(-11, -3, -45, -48, -34, -56, -30, -16, -6, 0)
    19 = @(120, 164, 312, 35, 180, 300, 92, 86, 211, 30)
    20 = @(648, 1811, 481, 1725, 1067, 2121, 890, 5031, 2838, 3744)
    21 = @("-", -21, "-", -125, -119, -96, -113, -740, -316, -304)
    22 = @(648, 1790, 481, 1600, 948, 2026, 777, 4291, 2522, 3439)
    23 = @("-", "-", "-", "-", "-", "-", "-", "-", "-", "-")
    24 = @(648, 1790, 481, 1600, 948, 2026, 777, 4291, 2522, 3439)
    25 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    26 = @(2072, 7476, 7827, 6996, 6413, 6629, 6214, 5887, 5256, 4017)
    27 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
}

foreach ($r in $rowData.Keys) {
    $vals = $rowData[$r]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range($cols[$i] + $r).Value = $vals[$i]
    }
}
